# Apply targeted cell value updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("M10").Value = 4
$ws.Range("N10").Value = 4

$ws.Range("N15").Value = 1

$ws.Range("M31").Value = 3
$ws.Range("N31").Value = 3

$ws.Range("M32").Value = 1
$ws.Range("N32").Value = 2

$ws.Range("M33").Value = 3
$ws.Range("N33").Value = 3

$ws.Range("M36").Value = 0
$ws.Range("N36").Value = 1

$ws.Range("M42").Value = 4
$ws.Range("N42").Value = 5

$ws.Range("M49").Value = 3
$ws.Range("N49").Value = 3
